$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before current row 2, shifting "Rerun new data",
# "Fix Shiny" and "Regional structure" down by one.
$ws.Rows.Item(2).Insert()

# New note next to "Rerun new data" (now row 3)
$ws.Range("B3").Value = "new frq file from Thom 4 July 12:09"

# New row 2 content
$ws.Range("A2").Value = "Explore Init"
$ws.Range("B2").Value = "try changing only M or only K"

# Autofit column B to best fit its new contents
$ws.Columns.Item(2).AutoFit()
